$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(44)
$tr = $sh.TextFrame.TextRange

# Apply edits back-to-front so earlier character offsets stay valid after
# each replacement changes the overall text length.

# Change 2 (last paragraph, originally starts at char 139): re-assert the
# paragraph's text ("*테이블 만드는 순서 : ...") over its full span so
# PowerPoint collapses the many single-word runs that made it up into one
# run.
$last = $tr.Characters(139, 28)
$last.Text = "*테이블 만드는 순서 : 회원>게시판>첨부파일>댓글"

# Change 1: inside the first bullet's "FK(게시물코드)" run, retype just the
# inner "게시물코드" as "회원코드" (position 36, length 5 chars) - this
# reproduces the three-way run split (FK( / 회원코드 / )) that PowerPoint
# produces when you select-and-retype the middle of a run.
$inner = $tr.Characters(36, 5)
$inner.Text = "회원코드"

# The text edits above make the autofit textbox relayout (its stored
# height drifts a hair from the saved original). Re-pin the box back to
# its on-disk size in points (112.7pt tall, 394.65pt wide) so it round-trips
# to the same EMU extent PowerPoint itself saved after this edit.
$sh.Width = 394.65001
$sh.Height = 112.70001
